# Weekly update: insert two new price-record rows (new rows 117 and 118)
# above the existing "Arveja Verde" records, pushing the old rows 117-138
# down to 119-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 117 (shifts old 117..138 -> 119..140)
$ws.Range("A117:A118").EntireRow.Insert()

# ---- New row 117 ----
$ws.Cells.Item(117, 1).Value  = 6
$ws.Cells.Item(117, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(117, 3).Value  = "Metropolitana"
$ws.Cells.Item(117, 4).Value  = 44476
$ws.Cells.Item(117, 5).Value  = 13
$ws.Cells.Item(117, 6).Value  = 100112022
$ws.Cells.Item(117, 7).Value  = "Arveja Verde"
$ws.Cells.Item(117, 8).Value  = "Perfection"
$ws.Cells.Item(117, 9).Value  = "Primera"
$ws.Cells.Item(117, 10).Value = 270
$ws.Cells.Item(117, 11).Value = 23000
$ws.Cells.Item(117, 12).Value = 24000
$ws.Cells.Item(117, 13).Value = 23630
$ws.Cells.Item(117, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(117, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(117, 16).Value = 945
$ws.Cells.Item(117, 17).Value = 25
$ws.Cells.Item(117, 18).Value = "Hortaliza"

# ---- New row 118 ----
$ws.Cells.Item(118, 1).Value  = 6
$ws.Cells.Item(118, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(118, 3).Value  = "Metropolitana"
$ws.Cells.Item(118, 4).Value  = 44476
$ws.Cells.Item(118, 5).Value  = 13
$ws.Cells.Item(118, 6).Value  = 100112022
$ws.Cells.Item(118, 7).Value  = "Arveja Verde"
$ws.Cells.Item(118, 8).Value  = "Sin especificar"
$ws.Cells.Item(118, 9).Value  = "Primera"
$ws.Cells.Item(118, 10).Value = 130
$ws.Cells.Item(118, 11).Value = 19000
$ws.Cells.Item(118, 12).Value = 20000
$ws.Cells.Item(118, 13).Value = 19385
$ws.Cells.Item(118, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(118, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(118, 16).Value = 775
$ws.Cells.Item(118, 17).Value = 25
$ws.Cells.Item(118, 18).Value = "Hortaliza"
